# "Real results" sheet: fill in measurements for crank positions 28-34.
#
# NOTE: this COM-interop runtime does not pre-populate the usual named
# Excel enum variables (e.g. $xlShiftDown, $xlPasteFormats all resolve to
# empty/null), so the well-known literal integer values are used instead:
#   xlShiftDown            = -4121
#   xlFormatFromLeftOrAbove =    0
#   xlPasteFormats         = -4122

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Real results")

# Insert 3 new rows above row 18 (old row 18 -> row 21), copying the format
# of the row directly above the insertion point so borders/number formats
# carry over just like a manual "insert copied cells" in the UI would.
$ws.Rows.Item(18).Resize(3).EntireRow.Insert(-4121, 0)

$ws.Range("B17:G17").Copy()
$ws.Range("B18:G18").PasteSpecial(-4122)
$ws.Range("B17:G17").Copy()
$ws.Range("B19:G19").PasteSpecial(-4122)
$ws.Range("B17:G17").Copy()
$ws.Range("B20:G20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the previously-empty rows 15-17 with the real measurements.
$ws.Cells.Item(15, 3).Value = 28
$ws.Cells.Item(15, 4).Value = 79
$ws.Cells.Item(15, 5).Value = 116
$ws.Cells.Item(15, 6).Formula = "=1/1.432"
$ws.Cells.Item(15, 7).Value = 2.92

$ws.Cells.Item(16, 3).Value = 29
$ws.Cells.Item(16, 4).Value = 101
$ws.Cells.Item(16, 5).Value = 135
$ws.Cells.Item(16, 6).Formula = "=1/1.383"
$ws.Cells.Item(16, 7).Value = 2.76

$ws.Cells.Item(17, 3).Value = 30
$ws.Cells.Item(17, 4).Value = 101
$ws.Cells.Item(17, 5).Value = 112
$ws.Cells.Item(17, 6).Formula = "=1/1.333"
$ws.Cells.Item(17, 7).Value = 2.56

# New row 18 (crank position 31). D18 gets its own one-off number-format
# tweak (still "0.0" but Excel records a distinct style id for it).
$ws.Cells.Item(18, 2).Value = "B"
$ws.Cells.Item(18, 3).Value = 31
$ws.Cells.Item(18, 4).Value = 107
$ws.Cells.Item(18, 4).NumberFormat = "0.0"
$ws.Cells.Item(18, 5).Value = 97
$ws.Cells.Item(18, 6).Formula = "=1/1.294"
$ws.Cells.Item(18, 7).Value = 2.4

# New row 19 (crank position 32) - the tape measurements were estimated,
# recorded as text with a trailing asterisk instead of plain numbers.
$ws.Cells.Item(19, 2).Value = "B"
$ws.Cells.Item(19, 3).Value = 32
$ws.Cells.Item(19, 4).Value = "119*"
$ws.Cells.Item(19, 5).Value = "116*"
$ws.Cells.Item(19, 6).Formula = "=1/1.258"
$ws.Cells.Item(19, 7).Value = 2.32

# New row 20 (crank position 33).
$ws.Cells.Item(20, 2).Value = "B"
$ws.Cells.Item(20, 3).Value = 33
$ws.Cells.Item(20, 4).Value = 88
$ws.Cells.Item(20, 5).Value = 94
$ws.Cells.Item(20, 6).Formula = "=1/1.217"
$ws.Cells.Item(20, 7).Value = 2.11
